$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EINHELL")
$ws.Activate()
$ws.Range("J2").Value = "B3423"
[void]$ws.Range("J11").Select()
